# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first sheet) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3328
$ws1.Range("F3").Value = 744
$ws1.Range("F4").Value = 127
$ws1.Range("F5").Value = 6936
$ws1.Range("F6").Value = 2300
$ws1.Range("F8").Value = 92
$ws1.Range("F12").Value = 30
$ws1.Range("F13").Value = 166
$ws1.Range("F14").Value = 478
$ws1.Range("F16").Value = 131

# --- Sheet "全部类型" (fourth sheet) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3328
$ws4.Range("F4").Value = 744
$ws4.Range("F5").Value = 127
$ws4.Range("F6").Value = 6936
$ws4.Range("F7").Value = 2300
$ws4.Range("F9").Value = 92
$ws4.Range("F13").Value = 30
$ws4.Range("F14").Value = 166
$ws4.Range("F15").Value = 478
$ws4.Range("F17").Value = 131
